$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# so Excel does not auto-convert them (the source data stores these as text).
$textForceCells = @("D5", "D6", "D14", "D19", "D21", "D23", "D26", "D27", "D28", "D31", "D36", "D37", "D38", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.085.65"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "2.421.15"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "568.97"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "139.46"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.405.74"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "26.07"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").Value = "2.854.76"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "61.013.09"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "2.403.80"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "7.85"
$ws.Range("E19").Value = "  +8.39%  "
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "322.86"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "6.09"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").Value = "64.64"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").Value = "582.46"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").Value = "8.25"
$ws.Range("E28").Value = "  -9.39%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "1.42"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "152.01"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").Value = "41.12"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -7.08%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0277"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "143.13"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "3.51"
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("D49").Value = "0.587"
$ws.Range("D50").Value = "19.49"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  -3.43%  "

# Restore default styling on the forced cells so no stray number format remains applied.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
